$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(92, 8).Value = 273.7143
$ws.Cells.Item(92, 9).Value = 225.8125
$ws.Cells.Item(92, 10).Value = 427
$ws.Cells.Item(92, 11).Value = 225.8125
$ws.Cells.Item(92, 12).Value = 427
$ws.Cells.Item(92, 13).Value = 1022.1875
$ws.Cells.Item(92, 14).Value = -2923
$ws.Cells.Item(128, 8).Value = 59800
$ws.Cells.Item(128, 10).Value = 59800
$ws.Cells.Item(128, 12).Value = 59800
$ws.Cells.Item(128, 14).Value = -69760
$ws.Cells.Item(132, 8).Value = 20844802
$ws.Cells.Item(132, 9).Value = 20844802
$ws.Cells.Item(132, 10).Value = 0
$ws.Cells.Item(132, 11).Value = 62534406
$ws.Cells.Item(132, 12).Value = 0
$ws.Cells.Item(132, 13).ClearContents()
$ws.Cells.Item(132, 14).Value = -62531876
$ws.Cells.Item(137, 8).Value = 2830.3076
$ws.Cells.Item(137, 9).Value = 1716.8334
$ws.Cells.Item(137, 10).Value = 3784.7144
$ws.Cells.Item(137, 11).Value = 5150.5002
$ws.Cells.Item(137, 12).Value = 11354.1432
$ws.Cells.Item(137, 13).Value = -2600.5002
$ws.Cells.Item(137, 14).Value = -16454.1432

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 1430.9048
$ws.Cells.Item(2, 10).Value = 1914.5
$ws.Cells.Item(2, 12).Value = 1914.5
$ws.Cells.Item(2, 14).Value = -2140.5
$ws.Cells.Item(4, 8).Value = 683.3333
$ws.Cells.Item(4, 9).Value = 300
$ws.Cells.Item(4, 11).Value = 300
$ws.Cells.Item(4, 13).Value = -184
$ws.Cells.Item(6, 8).Value = 12500
$ws.Cells.Item(6, 9).Value = 0
$ws.Cells.Item(6, 10).Value = 12500
$ws.Cells.Item(6, 11).Value = 0
$ws.Cells.Item(6, 12).ClearContents()
$ws.Cells.Item(6, 13).Value = 12500
$ws.Cells.Item(6, 14).Value = -12846
$ws.Cells.Item(9, 8).Value = 0
$ws.Cells.Item(9, 10).Value = 0
$ws.Cells.Item(9, 12).ClearContents()
$ws.Cells.Item(9, 14).Value = 0
$ws.Cells.Item(20, 8).Value = 0
$ws.Cells.Item(20, 10).Value = 0
$ws.Cells.Item(20, 12).ClearContents()
$ws.Cells.Item(20, 14).Value = 0
$ws.Cells.Item(23, 8).Value = 67659
$ws.Cells.Item(23, 10).Value = 67267.836
$ws.Cells.Item(23, 12).Value = 67267.836
$ws.Cells.Item(23, 14).Value = -67785.836
$ws.Cells.Item(37, 8).Value = 2900
$ws.Cells.Item(37, 9).Value = 2900
$ws.Cells.Item(37, 10).Value = 0
$ws.Cells.Item(37, 11).Value = 2900
$ws.Cells.Item(37, 12).ClearContents()
$ws.Cells.Item(37, 13).Value = -2627
$ws.Cells.Item(37, 14).Value = 0
$ws.Cells.Item(44, 8).Value = 40714
$ws.Cells.Item(44, 10).Value = 40714
$ws.Cells.Item(44, 12).Value = 40714
$ws.Cells.Item(44, 14).Value = -41690
$ws.Cells.Item(55, 8).Value = 23053
$ws.Cells.Item(55, 9).Value = 0
$ws.Cells.Item(55, 11).Value = 0
$ws.Cells.Item(55, 13).ClearContents()
$ws.Cells.Item(61, 8).Value = 6252123.5
$ws.Cells.Item(61, 9).Value = 8334978.5
$ws.Cells.Item(61, 11).Value = 8334978.5
$ws.Cells.Item(61, 13).Value = -8334766.5
$ws.Cells.Item(74, 8).Value = 1352.4546
$ws.Cells.Item(74, 9).Value = 1346.2693
$ws.Cells.Item(74, 10).Value = 1375.4286
$ws.Cells.Item(74, 11).Value = 1346.2693
$ws.Cells.Item(74, 12).Value = 1375.4286
$ws.Cells.Item(74, 13).Value = -472.2692999999999
$ws.Cells.Item(74, 14).Value = -3123.4286
$ws.Cells.Item(77, 8).Value = 1352.4546
$ws.Cells.Item(77, 9).Value = 1346.2693
$ws.Cells.Item(77, 10).Value = 1375.4286
$ws.Cells.Item(77, 11).Value = 6731.3465
$ws.Cells.Item(77, 12).Value = 6877.143
$ws.Cells.Item(77, 13).Value = -2363.3465
$ws.Cells.Item(77, 14).Value = -15613.143
$ws.Cells.Item(80, 8).Value = 20000
$ws.Cells.Item(80, 10).Value = 20000
$ws.Cells.Item(80, 12).Value = 20000
$ws.Cells.Item(80, 14).Value = -21996
$ws.Cells.Item(83, 8).Value = 20000
$ws.Cells.Item(83, 10).Value = 20000
$ws.Cells.Item(83, 12).Value = 60000
$ws.Cells.Item(83, 14).Value = -69984
$ws.Cells.Item(116, 8).Value = 1430.9048
$ws.Cells.Item(116, 10).Value = 1914.5
$ws.Cells.Item(116, 12).Value = 1914.5
$ws.Cells.Item(116, 14).Value = -6502.5
$ws.Cells.Item(136, 8).Value = 6252123.5
$ws.Cells.Item(136, 9).Value = 8334978.5
$ws.Cells.Item(136, 11).Value = 25004935.5
$ws.Cells.Item(136, 13).Value = -25002385.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 1430.9048
$ws.Cells.Item(3, 10).Value = 1914.5
$ws.Cells.Item(3, 12).Value = 1914.5
$ws.Cells.Item(3, 14).Value = -2142.5
$ws.Cells.Item(94, 8).Value = 2008.8889
$ws.Cells.Item(94, 9).Value = 1796.6666
$ws.Cells.Item(94, 10).Value = 2433.3333
$ws.Cells.Item(94, 11).Value = 1796.6666
$ws.Cells.Item(94, 12).Value = 2433.3333
$ws.Cells.Item(94, 13).Value = -1345.6666
$ws.Cells.Item(94, 14).Value = -3335.3333
$ws.Cells.Item(107, 8).Value = 13159181
$ws.Cells.Item(107, 9).Value = 20834126
$ws.Cells.Item(107, 10).Value = 2132.5
$ws.Cells.Item(107, 11).Value = 20834126
$ws.Cells.Item(107, 12).Value = 2132.5
$ws.Cells.Item(107, 13).Value = -20832206
$ws.Cells.Item(107, 14).Value = -5972.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 291.9091
$ws.Cells.Item(7, 9).Value = 165.85715
$ws.Cells.Item(7, 10).Value = 512.5
$ws.Cells.Item(7, 11).Value = 165.85715
$ws.Cells.Item(7, 12).Value = 512.5
$ws.Cells.Item(7, 13).Value = -52.85714999999999
$ws.Cells.Item(7, 14).Value = -738.5
$ws.Cells.Item(22, 8).Value = 410.18182
$ws.Cells.Item(22, 9).Value = 445.77777
$ws.Cells.Item(22, 10).Value = 250
$ws.Cells.Item(22, 11).Value = 445.77777
$ws.Cells.Item(22, 12).Value = 250
$ws.Cells.Item(22, 13).Value = -95.77776999999998
$ws.Cells.Item(22, 14).Value = -950
$ws.Cells.Item(31, 8).Value = 1213.8
$ws.Cells.Item(31, 9).Value = 903.44446
$ws.Cells.Item(31, 10).Value = 4007
$ws.Cells.Item(31, 11).Value = 903.44446
$ws.Cells.Item(31, 12).Value = 4007
$ws.Cells.Item(31, 13).Value = -608.44446
$ws.Cells.Item(31, 14).Value = -4597
$ws.Cells.Item(34, 8).Value = 1213.8
$ws.Cells.Item(34, 9).Value = 903.44446
$ws.Cells.Item(34, 10).Value = 4007
$ws.Cells.Item(34, 11).Value = 903.44446
$ws.Cells.Item(34, 12).Value = 4007
$ws.Cells.Item(34, 13).Value = -701.44446
$ws.Cells.Item(34, 14).Value = -4411
$ws.Cells.Item(99, 8).Value = 166668000
$ws.Cells.Item(99, 9).Value = 333334000
$ws.Cells.Item(99, 10).Value = 2000
$ws.Cells.Item(99, 11).Value = 333334000
$ws.Cells.Item(99, 12).Value = 2000
$ws.Cells.Item(99, 13).Value = -333332502
$ws.Cells.Item(99, 14).Value = -4996
$ws.Cells.Item(126, 8).Value = 166668000
$ws.Cells.Item(126, 9).Value = 333334000
$ws.Cells.Item(126, 10).Value = 2000
$ws.Cells.Item(126, 11).Value = 1000002000
$ws.Cells.Item(126, 12).Value = 6000
$ws.Cells.Item(126, 13).Value = -999999530
$ws.Cells.Item(126, 14).Value = -10940
$ws.Cells.Item(127, 8).Value = 30000
$ws.Cells.Item(127, 10).Value = 30000
$ws.Cells.Item(127, 12).Value = 30000
$ws.Cells.Item(127, 14).Value = -39920

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 37638348
$ws.Cells.Item(5, 9).Value = 77778200
$ws.Cells.Item(5, 10).Value = 7233.5625
$ws.Cells.Item(5, 11).Value = 233334600
$ws.Cells.Item(5, 12).Value = 21700.6875
$ws.Cells.Item(5, 13).Value = -233334488
$ws.Cells.Item(5, 14).Value = -21924.6875
$ws.Cells.Item(131, 8).Value = 916.14
$ws.Cells.Item(131, 9).Value = 876.6667
$ws.Cells.Item(131, 10).Value = 917.3608400000001
$ws.Cells.Item(131, 11).Value = 2630.0001
$ws.Cells.Item(131, 12).Value = 2752.08252
$ws.Cells.Item(131, 13).Value = 2409.9999
$ws.Cells.Item(131, 14).Value = -12832.08252
$ws.Cells.Item(135, 8).Value = 37638348
$ws.Cells.Item(135, 9).Value = 77778200
$ws.Cells.Item(135, 10).Value = 7233.5625
$ws.Cells.Item(135, 11).Value = 700003800
$ws.Cells.Item(135, 12).Value = 65102.0625
$ws.Cells.Item(135, 13).Value = -700001265
$ws.Cells.Item(135, 14).Value = -70172.0625

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 94.64286
$ws.Cells.Item(2, 9).Value = 21.625
$ws.Cells.Item(2, 10).Value = 192
$ws.Cells.Item(2, 11).Value = 21.625
$ws.Cells.Item(2, 12).Value = 192
$ws.Cells.Item(2, 13).Value = 91.375
$ws.Cells.Item(2, 14).Value = -418
$ws.Cells.Item(107, 8).Value = 262.26923
$ws.Cells.Item(107, 9).Value = 177.5
$ws.Cells.Item(107, 10).Value = 361.16666
$ws.Cells.Item(107, 11).Value = 177.5
$ws.Cells.Item(107, 12).Value = 361.16666
$ws.Cells.Item(107, 13).Value = 1742.5
$ws.Cells.Item(107, 14).Value = -4201.16666
$ws.Cells.Item(126, 8).Value = 4547565.5
$ws.Cells.Item(126, 9).Value = 1595.4
$ws.Cells.Item(126, 10).Value = 8335874
$ws.Cells.Item(126, 11).Value = 4786.200000000001
$ws.Cells.Item(126, 12).Value = 25007622
$ws.Cells.Item(126, 13).Value = -2316.200000000001
$ws.Cells.Item(126, 14).Value = -25012562

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 1600
$ws.Cells.Item(7, 9).Value = 1333.3334
$ws.Cells.Item(7, 11).Value = 1333.3334
$ws.Cells.Item(7, 13).Value = -1221.3334
$ws.Cells.Item(40, 8).Value = 16669260
$ws.Cells.Item(40, 9).Value = 2254.3333
$ws.Cells.Item(40, 10).Value = 41669770
$ws.Cells.Item(40, 11).Value = 2254.3333
$ws.Cells.Item(40, 12).Value = 41669770
$ws.Cells.Item(40, 13).Value = -2118.3333
$ws.Cells.Item(40, 14).Value = -41670042
$ws.Cells.Item(93, 8).Value = 1328.091
$ws.Cells.Item(93, 9).Value = 1318.1666
$ws.Cells.Item(93, 10).Value = 1340
$ws.Cells.Item(93, 11).Value = 1318.1666
$ws.Cells.Item(93, 12).Value = 1340
$ws.Cells.Item(93, 13).Value = -70.16660000000002
$ws.Cells.Item(93, 14).Value = -3836
$ws.Cells.Item(126, 8).Value = 1600
$ws.Cells.Item(126, 9).Value = 1333.3334
$ws.Cells.Item(126, 11).Value = 4000.0002
$ws.Cells.Item(126, 13).Value = -1530.0002

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(100, 8).Value = 802.25
$ws.Cells.Item(100, 10).Value = 802.25
$ws.Cells.Item(100, 12).Value = 1604.5
$ws.Cells.Item(100, 14).Value = -2686.5
$ws.Cells.Item(126, 8).Value = 3557.5715
$ws.Cells.Item(126, 9).Value = 1968
$ws.Cells.Item(126, 10).Value = 4749.75
$ws.Cells.Item(126, 11).Value = 5904
$ws.Cells.Item(126, 12).Value = 14249.25
$ws.Cells.Item(126, 13).Value = -3434
$ws.Cells.Item(126, 14).Value = -19189.25
